$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Standard CSRp"

$r = $ws2.Range("A3:I3")
$r.Merge()
$r.Value = "Dense matrix"
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.Borders.Weight = 2
$r.Borders.ColorIndex = 64

Write-Host "done"
